$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old secondary header row (old row 2); this shifts rows 3..32 up to 2..31
$ws.Rows.Item(2).Delete()

# Rewrite row 1 as a full header row
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Update selection to match target
$ws.Range("A2:K2").Select()
